$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.233.32'
$ws.Range("E2").Value = '  +1.71%  '
$ws.Range("D3").Value = '1.645.14'
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.86'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.56%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.506'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.47%  '
$ws.Range("E7").Value = '  -0.17%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("E9").Value = '  -0.28%  '
$ws.Range("E10").Value = '  +1.16%  '
$ws.Range("E11").Value = '  +0.23%  '
$ws.Range("E12").Value = '  +0.30%  '
$ws.Range("D13").Value = '1.873.00'
$ws.Range("E13").Value = '  +0.51%  '
$ws.Range("D14").Value = '1.644.19'
$ws.Range("E14").Value = '  +0.42%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.550'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.07%  '
$ws.Range("D16").Value = '0.0₃0764'
$ws.Range("E16").Value = '  -0.66%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.50'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.49%  '
$ws.Range("D18").Value = '26.209.38'
$ws.Range("E18").Value = '  +1.50%  '
$ws.Range("E19").Value = '  -0.12%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '195.40'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.25%  '
$ws.Range("E21").Value = '  -0.85%  '
$ws.Range("E22").Value = '  +0.55%  '
$ws.Range("E23").Value = '  -0.69%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '143.34'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.75%  '
$ws.Range("E25").Value = '  -0.15%  '
$ws.Range("E26").Value = '  -2.03%  '
$ws.Range("E27").Value = '  +1.76%  '
$ws.Range("E28").Value = '  -0.37%  '
$ws.Range("E29").Value = '  +0.39%  '
$ws.Range("E30").Value = '  +1.24%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0503'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.96%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.34'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.02%  '
$ws.Range("E34").Value = '  +1.32%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.40'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.89%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.914'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.88%  '
$ws.Range("D37").Value = '1.134.81'
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("E38").Value = '  +1.37%  '
$ws.Range("E39").Value = '  -1.53%  '
$ws.Range("E40").Value = '  +0.76%  '
$ws.Range("E42").Value = '  +1.70%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.00'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.35%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.798'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.98%  '
$ws.Range("D45").Value = '1.782.27'
$ws.Range("E45").Value = '  +0.52%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '56.16'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.59%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.48'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.40%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0517'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.80%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.71'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.64%  '
$ws.Range("E51").Value = '  +1.69%  '
